# Generate Report for Handback
#
# Updates the handback timestamps that were generated for the
# "2f7780e8-e76e-469e-91b9-bc99e9da7c12.md" file (row 3 on each sheet):
#   - Overview!G3            -> Latest HO Xliff Generate Date
#   - zh-cn!H3 / zh-cn!K3     -> Correspond Handoff / Handback Datetime
#   - de-de!H3 / de-de!K3     -> Correspond Handoff / Handback Datetime

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-05 20:57:08"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-05 20:56:58"
$wsZhCn.Range("K3").Value = "2016-09-05 20:57:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-05 20:57:08"
$wsDeDe.Range("K3").Value = "2016-09-05 20:57:35"
